# Insert 3 new data rows at the top of the price-record block (row 137),
# pushing the existing records (old rows 137-230) down to rows 140-233.
# Excel's Insert() inherits formatting from the row above, so column D
# (dates) keeps its existing date style/number-format automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("137:139").Insert()

# ---- New row 137 ----
$ws.Range("A137").Value = 11
$ws.Range("B137").Value = "Vega Monumental Concepción"
$ws.Range("C137").Value = "Bíobío"
$ws.Range("D137").Value = 44673
$ws.Range("E137").Value = 8
$ws.Range("F137").Value = 100114001
$ws.Range("G137").Value = "Papa"
$ws.Range("H137").Value = "Asterix"
$ws.Range("I137").Value = "1a (cosecha lavada)"
$ws.Range("J137").Value = 350
$ws.Range("K137").Value = 8500
$ws.Range("L137").Value = 9000
$ws.Range("M137").Value = 8786
$ws.Range("N137").Value = "$/malla 25 kilos"
$ws.Range("O137").Value = "Región de La Araucanía"
$ws.Range("P137").Value = 351
$ws.Range("Q137").Value = 25
$ws.Range("R137").Value = "Hortaliza"

# ---- New row 138 ----
$ws.Range("A138").Value = 11
$ws.Range("B138").Value = "Vega Monumental Concepción"
$ws.Range("C138").Value = "Bíobío"
$ws.Range("D138").Value = 44673
$ws.Range("E138").Value = 8
$ws.Range("F138").Value = 100114001
$ws.Range("G138").Value = "Papa"
$ws.Range("H138").Value = "Asterix"
$ws.Range("I138").Value = "1a (cosecha)"
$ws.Range("J138").Value = 350
$ws.Range("K138").Value = 7000
$ws.Range("L138").Value = 7500
$ws.Range("M138").Value = 7214
$ws.Range("N138").Value = "$/saco 25 kilos"
$ws.Range("O138").Value = "Región de Los Lagos"
$ws.Range("P138").Value = 289
$ws.Range("Q138").Value = 25
$ws.Range("R138").Value = "Hortaliza"

# ---- New row 139 ----
$ws.Range("A139").Value = 11
$ws.Range("B139").Value = "Vega Monumental Concepción"
$ws.Range("C139").Value = "Bíobío"
$ws.Range("D139").Value = 44673
$ws.Range("E139").Value = 8
$ws.Range("F139").Value = 100114001
$ws.Range("G139").Value = "Papa"
$ws.Range("H139").Value = "Patagonia"
$ws.Range("I139").Value = "1a (cosecha)"
$ws.Range("J139").Value = 350
$ws.Range("K139").Value = 6000
$ws.Range("L139").Value = 7000
$ws.Range("M139").Value = 6429
$ws.Range("N139").Value = "$/saco 25 kilos"
$ws.Range("O139").Value = "Región de Los Lagos"
$ws.Range("P139").Value = 257
$ws.Range("Q139").Value = 25
$ws.Range("R139").Value = "Hortaliza"
